$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Sheet 1": append three more "BC" / "IND2" rows (20, 21, 22) that mirror
#    the existing row 19 formatting/values, but with a blank "value" (column
#    C) -- this models rows whose indicator value came back NULL from the
#    importer (see commit message: "Allow null value on importer").
# ---------------------------------------------------------------------------
$data = $wb.Worksheets.Item("Sheet 1")

$srcRow = 19
$srcRange = $data.Range("A$srcRow`:H$srcRow")
foreach ($r in 20, 21, 22) {
    $dstRange = $data.Range("A$r`:H$r")
    $srcRange.Copy($dstRange)
    # Column C ("value") comes back blank/NULL for these records.
    $data.Range("C$r").Value2 = $null
}

# ---------------------------------------------------------------------------
# 2) Fix the "key has a dot on it" / blank-value query bug: the COUNT sheets
#    should not count rows where 'Sheet 1'!C:C (the value column) is blank.
#    Append the extra COUNTIFS criteria to both the per-geom COUNT sheet and
#    the "COUNT (Upper Level)" sheet.
# ---------------------------------------------------------------------------
$count = $wb.Worksheets.Item("COUNT")
for ($r = 2; $r -le 7; $r++) {
    $old = $count.Range("D$r").Formula
    $new = $old.Substring(0, $old.Length - 1) + ",'Sheet 1'!C:C,""<>""&""""" + ")"
    $count.Range("D$r").Formula = $new
}

$countUpper = $wb.Worksheets.Item("COUNT (Upper Level)")
for ($r = 2; $r -le 4; $r++) {
    $old = $countUpper.Range("D$r").Formula
    $new = $old.Substring(0, $old.Length - 1) + ",'Sheet 1'!C:C,""<>""&""""" + ")"
    $countUpper.Range("D$r").Formula = $new
}
